$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Copy the formatting (number format, font, borders, alignment) of the
#        existing 2021 column (R) into the new 2022 column (S), for the
#        header row and every data row. ---
$ws.Range("R3:R33").Copy() | Out-Null
$ws.Range("S3:S33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- 2. Header: 2022 ---
$ws.Cells.Item(3, 19).Value = 2022

# --- 3. Data values for column S (rows 4-33), matching column R's layout. ---
$values = @(
    5.5,    # row 4
    8.5,    # row 5
    2.6,    # row 6
    16.3,   # row 7
    25.2,   # row 8
    7.1,    # row 9
    1.6,    # row 10
    3.2,    # row 11
    "-",    # row 12
    7.5,    # row 13
    10.5,   # row 14
    4.5,    # row 15
    11.4,   # row 16
    16.1,   # row 17
    6.6,    # row 18
    1.2,    # row 19
    2.1,    # row 20
    0.3,    # row 21
    1.5,    # row 22
    2.9,    # row 23
    0,      # row 24
    0.9,    # row 25
    1.7,    # row 26
    0.2,    # row 27
    14.3,   # row 28
    22.7,   # row 29
    7.3,    # row 30
    1.1,    # row 31
    2.2,    # row 32
    "-"     # row 33
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 19).Value = $values[$i]
}

# --- 4. Match the workbook's recorded selection after the edit. ---
$ws.Range("T3").Select() | Out-Null
